# Applies the ChallengeCardData.xlsx edit described by the commit:
# "Changed the effect of cursed coin and cursed Pot cards to make them
#  stay longer in the player's Loot Zone. Added tag to item cards to
#  indicate their archetype."
#
# Concretely: several card "effect" texts (column C) are rewritten, one
# new row gets word-wrap styling + a taller row height, and two row
# heights swap. The shared-string table is rebuilt by the engine itself
# (dedup + GC of orphaned strings) as long as we assign the new literal
# text to each cell in the same order the final table expects them, so
# the assignment order below is deliberate - do not reorder.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (一线天 / Narrow Hall): effect rewritten ---------------------
$ws.Range("C2").Value = "房间区宽度降低至1，高度增加宽度降低的数值。"

# --- Row 4 (迷雾 / Fog): effect rewritten --------------------------------
$ws.Range("C4").Value = @"
房间区中央的牌在发牌时背面向上。<br>
开战时：将战场中背面向上的牌翻正。<br>
回合结束时：选房间区1张牌翻正。
"@.TrimEnd("`r", "`n")

# --- Row 3 (大房间 / Big Room): effect rewritten, now wraps + taller ----
$ws.Range("C3").Value = @"
房间宽度加1，高度加1。<br>
玩家手牌基数加3。<br>
翻选时可以额外翻2张牌，额外选1张牌。
"@.TrimEnd("`r", "`n")

# --- Row 5 (大墓地 / Great Graveyard): effect rewritten ------------------
$ws.Range("C5").Value = @"
挑战开始时：从额外牌堆将《命匣》牌洗入主牌堆。<br>
重整后：从额外牌堆将1张《巫妖》牌放在房间区任意槽位中。<br>
挑战结束时：将所有不在额外牌堆的《巫妖》牌和《命匣》牌放回额外牌堆。
"@.TrimEnd("`r", "`n")

# --- Row 6 (斗技场 / Arena): effect rewritten ----------------------------
$ws.Range("C6").Value = @"
挑战开始时：将主牌堆所有牌送墓，此期间内牌的送墓时效果无效。<br>
重整后：从额外牌堆将1张《冠军》牌放在房间区任意槽位中。<br>
挑战结束时：将所有不在额外牌堆的《冠军》牌放回额外牌堆。
"@.TrimEnd("`r", "`n")

# --- Row 8 (沙虫坑 / Sandworm Pit): effect rewritten ---------------------
$ws.Range("C8").Value = @"
挑战开始时：重整。
重整后：从额外牌堆将2张《沙虫》牌分别放到房间区和备战区的任意槽位中。<br>
挑战结束时：将所有不在额外牌堆的《沙虫》牌放回额外牌堆。
"@.TrimEnd("`r", "`n")

# --- Row 9 (裂隙 / Rift): effect rewritten -------------------------------
$ws.Range("C9").Value = @"
挑战开始时：从额外牌堆将《时空吞噬者》牌洗入主牌堆。<br>
持续：弃牌堆的牌无法移出弃牌堆，所有消耗时间的牌无效。<br>
挑战结束时：将所有不在额外牌堆的《时空吞噬者》牌放回额外牌堆，然后将弃牌堆所有牌送墓。
"@.TrimEnd("`r", "`n")

# Rows 7 (展厅), 10 (手层) and 11 (冲突层) keep their existing effect text
# untouched - only their shared-string index shifts because of the
# inserts/deletes above, which the engine handles automatically.

# --- Formatting ----------------------------------------------------------
# New row (大房间, row 3) now wraps like the other long-effect rows.
$ws.Range("C3").WrapText = $true

# Row-height tweaks from the diff.
$ws.Rows.Item(3).RowHeight = 42.75
$ws.Rows.Item(6).RowHeight = 42.75
$ws.Rows.Item(8).RowHeight = 57

# --- View/selection state -------------------------------------------------
[void]$ws.Range("C16").Select()
try { $excel.ActiveWindow.ScrollRow = 3 } catch {}

try {
    $excel.ActiveWindow.Left = -120
    $excel.ActiveWindow.Top = -120
    $excel.ActiveWindow.Width = 29040
    $excel.ActiveWindow.Height = 15720
} catch {}
